# Updates FlashScore odds data (Jogos_da_Semana_FlashScore_2025-02-24.xlsx).
# Refreshes betting-odds columns (G:AS) for the matches below with updated quotes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: AS Roma vs Monza
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 4.33
$ws.Range("Q2").Value = 1.73
$ws.Range("R2").Value = 2.1
$ws.Range("S2").Value = 2.75
$ws.Range("T2").Value = 1.44
$ws.Range("U2").Value = 1.33
$ws.Range("V2").Value = 3.25
$ws.Range("W2").Value = 2.25
$ws.Range("X2").Value = 1.57
$ws.Range("Y2").Value = 6.5
$ws.Range("AF2").Value = 11
$ws.Range("AJ2").Value = 21
$ws.Range("AL2").Value = 29
$ws.Range("AO2").Value = 67

# Row 3: Sevilla vs Mallorca
$ws.Range("G3").Value = 1.91
$ws.Range("I3").Value = 4.33
$ws.Range("L3").Value = 5
$ws.Range("S3").Value = 4.5
$ws.Range("T3").Value = 1.2
$ws.Range("U3").Value = 1.53
$ws.Range("V3").Value = 2.38
$ws.Range("AG3").Value = 21
$ws.Range("AJ3").Value = 9
$ws.Range("AK3").Value = 21
$ws.Range("AM3").Value = 51
$ws.Range("AP3").Value = 1.86
$ws.Range("AQ3").Value = 2.04

# Row 4: Barracas Central vs Newells Old Boys
$ws.Range("H4").Value = 2.75
$ws.Range("I4").Value = 2.88
$ws.Range("K4").Value = 1.8
$ws.Range("O4").Value = 1.67
$ws.Range("P4").Value = 2.1
$ws.Range("Q4").Value = 3.4
$ws.Range("R4").Value = 1.33
$ws.Range("S4").Value = 7
$ws.Range("T4").Value = 1.1
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 2.08
$ws.Range("W4").Value = 2.38
$ws.Range("X4").Value = 1.53
$ws.Range("Y4").Value = 6
$ws.Range("AC4").Value = 34
$ws.Range("AE4").Value = 5
$ws.Range("AN4").Value = 34
$ws.Range("AO4").Value = 51
$ws.Range("AR4").Value = 5.6
$ws.Range("AS4").Value = 1.15

# Row 5: Argentinos Jrs vs Ind. Rivadavia
$ws.Range("G5").Value = 1.65
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 6
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("U5").Value = 1.53
$ws.Range("V5").Value = 2.38
$ws.Range("AA5").Value = 9
$ws.Range("AE5").Value = 7
$ws.Range("AJ5").Value = 11
$ws.Range("AK5").Value = 29
$ws.Range("AL5").Value = 21
$ws.Range("AP5").Value = 1.88
$ws.Range("AQ5").Value = 1.98

# Row 6: Tigre vs Banfield
$ws.Range("I6").Value = 4
$ws.Range("W6").Value = 2.5
$ws.Range("X6").Value = 1.5
$ws.Range("Z6").Value = 8
$ws.Range("AF6").Value = 6.5
$ws.Range("AK6").Value = 19
$ws.Range("AR6").Value = 5

# Row 7: Nueva Chicago vs Chacarita Juniors
$ws.Range("M7").Value = 1.17
$ws.Range("N7").Value = 5
$ws.Range("Q7").Value = 3.4
$ws.Range("R7").Value = 1.33
$ws.Range("S7").Value = 7
$ws.Range("T7").Value = 1.1
$ws.Range("V7").Value = 2.08

# Row 8: Racing Cordoba vs Alvarado
$ws.Range("Q8").Value = 3.1
$ws.Range("R8").Value = 1.36

# Row 11: Colo Colo vs O'Higgins
$ws.Range("G11").Value = 1.45
$ws.Range("H11").Value = 4.1
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 13
$ws.Range("Q11").Value = 1.8
$ws.Range("R11").Value = 2
$ws.Range("AG11").Value = 19
$ws.Range("AI11").Value = 351
$ws.Range("AP11").Value = 1.38
$ws.Range("AQ11").Value = 3
$ws.Range("AR11").Value = 2.34
$ws.Range("AS11").Value = 1.58

# Row 17: Tecnico U. vs Libertad
$ws.Range("G17").Value = 1.67
$ws.Range("H17").Value = 3.7
$ws.Range("I17").Value = 5
$ws.Range("O17").Value = 1.36
$ws.Range("P17").Value = 3
$ws.Range("Q17").Value = 2.1
$ws.Range("R17").Value = 1.7
$ws.Range("Z17").Value = 7.5
$ws.Range("AK17").Value = 23
$ws.Range("AO17").Value = 41

# Row 18: Sheffield Utd vs Leeds
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 3.75
$ws.Range("J18").Value = 5
$ws.Range("M18").Value = 1.06
$ws.Range("N18").Value = 10
$ws.Range("Y18").Value = 10
$ws.Range("AA18").Value = 13
$ws.Range("AC18").Value = 34
$ws.Range("AE18").Value = 9.5
$ws.Range("AF18").Value = 7
$ws.Range("AI18").Value = 401

# Row 23: Helmond vs FC Emmen
$ws.Range("G23").Value = 2.8
$ws.Range("H23").Value = 3.6
$ws.Range("I23").Value = 2.3
$ws.Range("J23").Value = 3.25
$ws.Range("K23").Value = 2.3
$ws.Range("Q23").Value = 1.62
$ws.Range("R23").Value = 2.25
$ws.Range("U23").Value = 1.3
$ws.Range("V23").Value = 3.4
$ws.Range("W23").Value = 1.53
$ws.Range("X23").Value = 2.38
$ws.Range("Y23").Value = 12
$ws.Range("AA23").Value = 11
$ws.Range("AC23").Value = 21
$ws.Range("AE23").Value = 15
$ws.Range("AF23").Value = 7
$ws.Range("AH23").Value = 41
$ws.Range("AI23").Value = 126
$ws.Range("AJ23").Value = 11
$ws.Range("AK23").Value = 13
$ws.Range("AR23").Value = 2.05
$ws.Range("AS23").Value = 1.8

# Row 24: Jong Ajax vs Cambuur
$ws.Range("G24").Value = 2.63
$ws.Range("I24").Value = 2.3
$ws.Range("J24").Value = 3.1
$ws.Range("L24").Value = 2.88
$ws.Range("W24").Value = 1.44
$ws.Range("X24").Value = 2.63
$ws.Range("AA24").Value = 10
$ws.Range("AC24").Value = 19
$ws.Range("AD24").Value = 21
$ws.Range("AK24").Value = 15
$ws.Range("AM24").Value = 23
$ws.Range("AN24").Value = 17

# Row 25: Jong PSV vs Dordrecht
$ws.Range("G25").Value = 3.25
$ws.Range("H25").Value = 4.2
$ws.Range("I25").Value = 1.95
$ws.Range("J25").Value = 3.5
$ws.Range("L25").Value = 2.4
$ws.Range("N25").Value = 23
$ws.Range("O25").Value = 1.1
$ws.Range("P25").Value = 7
$ws.Range("Q25").Value = 1.36
$ws.Range("R25").Value = 3.1
$ws.Range("S25").Value = 1.83
$ws.Range("T25").Value = 1.83
$ws.Range("W25").Value = 1.36
$ws.Range("X25").Value = 3
$ws.Range("Y25").Value = 19
$ws.Range("Z25").Value = 23
$ws.Range("AA25").Value = 12
$ws.Range("AC25").Value = 21
$ws.Range("AE25").Value = 26
$ws.Range("AF25").Value = 9.5
$ws.Range("AG25").Value = 11
$ws.Range("AM25").Value = 21
$ws.Range("AN25").Value = 13
$ws.Range("AO25").Value = 17

# Row 26: Sportivo Trinidense vs General Caballero JLM
$ws.Range("M26").Value = 1.07
$ws.Range("N26").Value = 9
$ws.Range("O26").Value = 1.36
$ws.Range("P26").Value = 3
$ws.Range("Q26").Value = 2.15
$ws.Range("R26").Value = 1.67
$ws.Range("S26").Value = 4
$ws.Range("T26").Value = 1.22

# Row 27: Binacional vs Cusco
$ws.Range("G27").Value = 2.25
$ws.Range("H27").Value = 3.3
$ws.Range("I27").Value = 3.1
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 2.1
$ws.Range("L27").Value = 3.75
$ws.Range("Z27").Value = 11
$ws.Range("AA27").Value = 9.5
$ws.Range("AB27").Value = 21
$ws.Range("AC27").Value = 19
$ws.Range("AD27").Value = 29
$ws.Range("AG27").Value = 13
$ws.Range("AK27").Value = 15
$ws.Range("AP27").Value = 1.48
$ws.Range("AQ27").Value = 2.7
$ws.Range("AR27").Value = 2.65
$ws.Range("AS27").Value = 1.48

# Row 28: Motor Lublin vs GKS Katowice
$ws.Range("N28").Value = 13
$ws.Range("Q28").Value = 1.83
$ws.Range("R28").Value = 2.03

# Row 29: S. Wola vs Warta Poznan
$ws.Range("O29").Value = 1.4
$ws.Range("P29").Value = 2.75
$ws.Range("Q29").Value = 2.25
$ws.Range("R29").Value = 1.62

# Row 30: FC Porto vs Vitoria Guimaraes
$ws.Range("G30").Value = 1.8
$ws.Range("H30").Value = 3.5
$ws.Range("I30").Value = 4.5
$ws.Range("J30").Value = 2.5
$ws.Range("L30").Value = 5
$ws.Range("M30").Value = 1.06
$ws.Range("N30").Value = 10
$ws.Range("Q30").Value = 2.1
$ws.Range("R30").Value = 1.7
$ws.Range("S30").Value = 3.75
$ws.Range("T30").Value = 1.25
$ws.Range("W30").Value = 2
$ws.Range("X30").Value = 1.75
$ws.Range("Z30").Value = 8
$ws.Range("AB30").Value = 15
$ws.Range("AD30").Value = 29
$ws.Range("AE30").Value = 8.5
$ws.Range("AF30").Value = 6.5
$ws.Range("AG30").Value = 17
$ws.Range("AI30").Value = 401
$ws.Range("AJ30").Value = 11
$ws.Range("AK30").Value = 23
$ws.Range("AL30").Value = 15
$ws.Range("AO30").Value = 41
